$d = $word.ActiveDocument

# The eight text blocks that get cyclically rotated around the document.
$A = 'Fornecer oportunidade de realização de treinamento profissional de Engenharia Ambiental em empresa ou instituição sob supervisão de docente do Departamento de Ciências Básicas e Ambientais da EEL. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional.'
$B = '4780627 - Ana Lucia Gabas Ferreira'
$C = 'Processo seletivo. Plano de trabalho específico. Realização do estágio. Relatório final.'
$D = 'Participação do aluno em processo seletivo de empresas, instituições de pesquisa ou no setor acadêmico. O estágio realizado sob a supervisão de docente designado pelo Departamento de Ciências Básicas e Ambientais da Escola de Engenharia de Lorena. O conteúdo será estabelecido no Plano de Trabalho entre o supervisor responsável pelo Estágio e o docente supervisor. Apresentação de relatório final sobre as atividades desenvolvidas no estágio.'
$E = 'Supervisão das atividades desenvolvidas pelo aluno durante o estágio.'
$F = 'A nota final será baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio.'
$G = 'Devido às características da disciplina, não será oferecida recuperação.'
$H = 'Não há.'

# Stage 1: replace each old value with a unique placeholder token so the
# later replacements (which reuse these same strings elsewhere) don't
# collide with each other while the rotation is applied.
$map = @(
    @{ old = $A; tok = '@@TOK_A@@' },
    @{ old = $B; tok = '@@TOK_B@@' },
    @{ old = $C; tok = '@@TOK_C@@' },
    @{ old = $D; tok = '@@TOK_D@@' },
    @{ old = $E; tok = '@@TOK_E@@' },
    @{ old = $F; tok = '@@TOK_F@@' },
    @{ old = $G; tok = '@@TOK_G@@' },
    @{ old = $H; tok = '@@TOK_H@@' }
)

foreach ($item in $map) {
    $d.Content.Find.Execute($item.old, $true, $false, $false, $false, $false, $true, 1, $false, $item.tok, 2)
}

# Stage 2: replace each placeholder token with the new value at that
# position, per the rotation described by the diff:
#   Objetivos (was A)            -> C
#   Docente (was B)               -> A
#   Programa resumido (was C)     -> D
#   Programa (was D)              -> E
#   Método (was E)                -> F
#   Critério (was F)              -> G
#   Norma de recuperação (was G)  -> H
#   Bibliografia (was H)          -> B
$map2 = @(
    @{ tok = '@@TOK_A@@'; new = $C },
    @{ tok = '@@TOK_B@@'; new = $A },
    @{ tok = '@@TOK_C@@'; new = $D },
    @{ tok = '@@TOK_D@@'; new = $E },
    @{ tok = '@@TOK_E@@'; new = $F },
    @{ tok = '@@TOK_F@@'; new = $G },
    @{ tok = '@@TOK_G@@'; new = $H },
    @{ tok = '@@TOK_H@@'; new = $B }
)

foreach ($item in $map2) {
    $d.Content.Find.Execute($item.tok, $true, $false, $false, $false, $false, $true, 1, $false, $item.new, 2)
}
